# Update "想去人数" (want-to-go count) values in F column for the
# "展览" sheet and the corresponding rows in the "全部类型" sheet.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 7672
$ws1.Range("F3").Value = 296
$ws1.Range("F4").Value = 32
$ws1.Range("F5").Value = 470
$ws1.Range("F6").Value = 4361
$ws1.Range("F7").Value = 330
$ws1.Range("F8").Value = 615
$ws1.Range("F10").Value = 689

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 7672
$ws4.Range("F4").Value = 296
$ws4.Range("F5").Value = 32
$ws4.Range("F6").Value = 470
$ws4.Range("F7").Value = 4361
$ws4.Range("F8").Value = 330
$ws4.Range("F9").Value = 615
$ws4.Range("F11").Value = 689
